# Scheduled market-data refresh: update cached price/profit figures across
# the per-job Leve tables (ALC, ARM, BSM, CRP, GSM, LTW, WVR). Values are
# static snapshots (no formulas in this workbook), so each changed cell is
# written directly with its refreshed figure. A couple of rows on the GSM
# sheet (70 and 73) now only have a single combined profit figure, so their
# old NQ/HQ split collapses into the M column and the N column is cleared.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 2777.3635
$ws.Range("I100").Value = 1379
$ws.Range("K100").Value = 1379
$ws.Range("M100").Value = -838

# Row 106
$ws.Range("H106").Value = 1364.5
$ws.Range("I106").Value = 1364.5
$ws.Range("K106").Value = 1364.5
$ws.Range("M106").Value = -733.5

# Row 112
$ws.Range("H112").Value = 5749315.5
$ws.Range("J112").Value = 5816063
$ws.Range("L112").Value = 17448189
$ws.Range("N112").Value = -17450405

# Row 135
$ws.Range("H135").Value = 1531.5
$ws.Range("I135").Value = 848.25
$ws.Range("K135").Value = 7634.25
$ws.Range("M135").Value = -5099.25

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 875
$ws.Range("I2").Value = 887.6
$ws.Range("K2").Value = 887.6
$ws.Range("M2").Value = -774.6

# Row 32
$ws.Range("H32").Value = 9690.157999999999
$ws.Range("I32").Value = 9339.611000000001
$ws.Range("K32").Value = 9339.611000000001
$ws.Range("M32").Value = -9052.611000000001

# Row 74
$ws.Range("H74").Value = 25670438
$ws.Range("I74").Value = 26345896
$ws.Range("J74").Value = 2993
$ws.Range("K74").Value = 26345896
$ws.Range("L74").Value = 2993
$ws.Range("M74").Value = -26345022
$ws.Range("N74").Value = -4741

# Row 77
$ws.Range("H77").Value = 25670438
$ws.Range("I77").Value = 26345896
$ws.Range("J77").Value = 2993
$ws.Range("K77").Value = 131729480
$ws.Range("L77").Value = 14965
$ws.Range("M77").Value = -131725112
$ws.Range("N77").Value = -23701

# Row 97
$ws.Range("H97").Value = 1073.3784
$ws.Range("I97").Value = 1102.8182
$ws.Range("J97").Value = 1030.2
$ws.Range("K97").Value = 1102.8182
$ws.Range("L97").Value = 1030.2
$ws.Range("M97").Value = -606.8181999999999
$ws.Range("N97").Value = -2022.2

# Row 102
$ws.Range("H102").Value = 226129.44
$ws.Range("I102").Value = 402795.8
$ws.Range("K102").Value = 402795.8
$ws.Range("M102").Value = -401173.8

# Row 116
$ws.Range("H116").Value = 875
$ws.Range("I116").Value = 887.6
$ws.Range("K116").Value = 887.6
$ws.Range("M116").Value = 1406.4

# Row 122
$ws.Range("H122").Value = 3560.818
$ws.Range("I122").Value = 2724.9092
$ws.Range("K122").Value = 8174.7276
$ws.Range("M122").Value = -5724.7276

# Row 132
$ws.Range("H132").Value = 32312648
$ws.Range("I132").Value = 11003.12
$ws.Range("K132").Value = 33009.36
$ws.Range("M132").Value = -30479.36

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 875
$ws.Range("I3").Value = 887.6
$ws.Range("K3").Value = 887.6
$ws.Range("M3").Value = -773.6

# Row 64
$ws.Range("H64").Value = 1617.3334
$ws.Range("J64").Value = 1540.8
$ws.Range("L64").Value = 1540.8
$ws.Range("N64").Value = -1990.8

# Row 67
$ws.Range("H67").Value = 1617.3334
$ws.Range("J67").Value = 1540.8
$ws.Range("L67").Value = 1540.8
$ws.Range("N67").Value = -3100.8

# Row 94
$ws.Range("H94").Value = 1879.6923
$ws.Range("I94").Value = 522.25
$ws.Range("J94").Value = 2483
$ws.Range("K94").Value = 522.25
$ws.Range("L94").Value = 2483
$ws.Range("M94").Value = -71.25
$ws.Range("N94").Value = -3385

# Row 107
$ws.Range("H107").Value = 3480.5715
$ws.Range("I107").Value = 2936.25
$ws.Range("K107").Value = 2936.25
$ws.Range("M107").Value = -1016.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3295.6345
$ws.Range("I31").Value = 2219.7693
$ws.Range("K31").Value = 2219.7693
$ws.Range("M31").Value = -1924.7693

# Row 34
$ws.Range("H34").Value = 3295.6345
$ws.Range("I34").Value = 2219.7693
$ws.Range("K34").Value = 2219.7693
$ws.Range("M34").Value = -2017.7693

# Row 132
$ws.Range("H132").Value = 51984.074
$ws.Range("I132").Value = 59907.06
$ws.Range("J132").Value = 5766.6665
$ws.Range("K132").Value = 179721.18
$ws.Range("L132").Value = 17299.9995
$ws.Range("M132").Value = -177191.18
$ws.Range("N132").Value = -22359.9995

# Row 134
$ws.Range("H134").Value = 1368.0769
$ws.Range("I134").Value = 1381.3043
$ws.Range("J134").Value = 1266.6666
$ws.Range("K134").Value = 4143.9129
$ws.Range("L134").Value = 3799.9998
$ws.Range("M134").Value = -1608.9129
$ws.Range("N134").Value = -8869.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 412303.2
$ws.Range("I70").Value = 412303.2
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 412303.2
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -412033.2
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 412303.2
$ws.Range("I73").Value = 412303.2
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 412303.2
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -411367.2
$ws.Range("N73").ClearContents()

# Row 126
$ws.Range("H126").Value = 5939.3125
$ws.Range("I126").Value = 5419.0835
$ws.Range("K126").Value = 16257.2505
$ws.Range("M126").Value = -13787.2505

# Row 132
$ws.Range("H132").Value = 2461.125
$ws.Range("I132").Value = 2130.3333
$ws.Range("J132").Value = 2659.6
$ws.Range("K132").Value = 6390.999899999999
$ws.Range("L132").Value = 7978.799999999999
$ws.Range("M132").Value = -3860.999899999999
$ws.Range("N132").Value = -13038.8

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2155.9534
$ws.Range("I61").Value = 1464.4615
$ws.Range("J61").Value = 3213.5293
$ws.Range("K61").Value = 1464.4615
$ws.Range("L61").Value = 3213.5293
$ws.Range("M61").Value = -1262.4615
$ws.Range("N61").Value = -3617.5293

# Row 113
$ws.Range("H113").Value = 2155.9534
$ws.Range("I113").Value = 1464.4615
$ws.Range("J113").Value = 3213.5293
$ws.Range("K113").Value = 1464.4615
$ws.Range("L113").Value = 3213.5293
$ws.Range("M113").Value = 705.5385000000001
$ws.Range("N113").Value = -7553.5293

# Row 136
$ws.Range("H136").Value = 2175.1538
$ws.Range("I136").Value = 1727.8
$ws.Range("J136").Value = 3666.3333
$ws.Range("K136").Value = 5183.4
$ws.Range("L136").Value = 10998.9999
$ws.Range("M136").Value = -2633.4
$ws.Range("N136").Value = -16098.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1806.762
$ws.Range("I81").Value = 1833.7273
$ws.Range("J81").Value = 1777.1
$ws.Range("K81").Value = 3667.4546
$ws.Range("L81").Value = 3554.2
$ws.Range("M81").Value = -2606.4546
$ws.Range("N81").Value = -5676.2

# Row 84
$ws.Range("H84").Value = 1806.762
$ws.Range("I84").Value = 1833.7273
$ws.Range("J84").Value = 1777.1
$ws.Range("K84").Value = 18337.273
$ws.Range("L84").Value = 17771
$ws.Range("M84").Value = -13033.273
$ws.Range("N84").Value = -28379

# Row 107
$ws.Range("H107").Value = 380.36365
$ws.Range("I107").Value = 403.77777
$ws.Range("J107").Value = 275
$ws.Range("K107").Value = 1211.33331
$ws.Range("L107").Value = 825
$ws.Range("M107").Value = 708.66669
$ws.Range("N107").Value = -4665

# Row 110
$ws.Range("H110").Value = 79991.664
$ws.Range("J110").Value = 79991.664
$ws.Range("L110").Value = 79991.664
$ws.Range("N110").Value = -88171.664
